$wb = $excel.ActiveWorkbook

# Map of row -> new value for column F ("想去人数")
$updates = @{
    23 = 4332
    26 = 1150
    29 = 690
    31 = 347
    33 = 180
}

# Both "展览" and "全部类型" sheets contain the same data table and both
# need the same updates applied.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
